$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new weekly snapshot row (row 8), mirroring the style of the
# existing date rows (numeric serial date with "yyyy-mm-dd" format).
$ws.Range("A8").NumberFormat = "yyyy-mm-dd"
$ws.Range("A8").Value = (Get-Date -Year 2026 -Month 1 -Day 4 -Hour 0 -Minute 0 -Second 0)

$ws.Range("B8").NumberFormat = "yyyy-mm-dd"
$ws.Range("B8").Value = (Get-Date -Year 2033 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0)
